# Edit script for 北京-漫展信息.xlsx
# Applies the gh-pages data refresh: updated "interested" counts (column F)
# across sheets, plus a newly-scraped exhibition row inserted into 展览 (sheet 1).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (1): insert the new event as row 42, shifting rows 42-44 down to 43-45 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(42).Insert()

# Copy formatting (bold/border style) of the index cell down from the row below so the
# new row 42 index cell (A42) matches the sheet's existing look.
$ws1.Range("A43").Copy()
$ws1.Range("A42").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# New row 42 content (newly scraped event). B42 is entered through a temporary Text
# number format (then restored) so the "2024.05.01"-style string is kept as literal
# text instead of being auto-parsed into a date serial.
$ws1.Range("A42").Value = 41
$ws1.Range("B42").NumberFormat = "@"
$ws1.Range("B42").Value = "2024.05.01"
$ws1.Range("C42").Copy()
$ws1.Range("B42").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false
$ws1.Range("C42").Value = "北京·IDO动漫游戏嘉年华45th同人创作大会"
$ws1.Range("D42").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws1.Range("E42").Value = "2024.05.01 09:30-05.03 17:00"
$ws1.Range("F42").Value = 0
$ws1.Range("G42").Value = 75
$ws1.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=82011"
$ws1.Range("I42").Value = "//i0.hdslb.com/bfs/openplatform/202402/2Aw7PvCg1708656416512.png"

# Rows 43, 44 and 45 (the old rows 42-44, now shifted down by the insert) keep their
# original content, but the sequential index in column A is recomputed for the new
# row position, and two of the three pick up an updated "interested" count.
$ws1.Range("A43").Value = 42
$ws1.Range("A44").Value = 43
$ws1.Range("A45").Value = 44
$ws1.Range("F43").Value = 1898
$ws1.Range("F44").Value = 2153

# --- Sheet "展览" (1): refresh "interested" counts (column F) for unaffected rows ---
$ws1.Range("F5").Value = 353
$ws1.Range("F6").Value = 407
$ws1.Range("F7").Value = 902
$ws1.Range("F9").Value = 553
$ws1.Range("F12").Value = 1176
$ws1.Range("F17").Value = 6731
$ws1.Range("F19").Value = 78
$ws1.Range("F21").Value = 7655
$ws1.Range("F24").Value = 3422
$ws1.Range("F25").Value = 35
$ws1.Range("F26").Value = 2156
$ws1.Range("F29").Value = 188
$ws1.Range("F34").Value = 204
$ws1.Range("F35").Value = 1787
$ws1.Range("F37").Value = 202
$ws1.Range("F39").Value = 9
$ws1.Range("F41").Value = 1255

# --- Sheet "演出" (2): refresh "interested" count (column F) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 7

# --- Sheet "全部类型" (4): refresh "interested" counts (column F) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F7").Value = 353
$ws4.Range("F8").Value = 407
$ws4.Range("F9").Value = 902
$ws4.Range("F11").Value = 553
$ws4.Range("F14").Value = 1176
$ws4.Range("F20").Value = 6732
$ws4.Range("F22").Value = 78
$ws4.Range("F24").Value = 7655
$ws4.Range("F27").Value = 3422
$ws4.Range("F28").Value = 35
$ws4.Range("F29").Value = 2156
$ws4.Range("F32").Value = 188
$ws4.Range("F38").Value = 1787
$ws4.Range("F40").Value = 202
$ws4.Range("F42").Value = 9
$ws4.Range("F44").Value = 1255
$ws4.Range("F45").Value = 1898
$ws4.Range("F47").Value = 2153
